$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 943.2366
$ws.Range("I15").Value = 943.2366
$ws.Range("K15").Value = 2829.7098
$ws.Range("M15").Value = -2660.7098

# row 32
$ws.Range("H32").Value = 540
$ws.Range("I32").Value = 633.3333
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 633.3333
$ws.Range("L32").Value = 400
$ws.Range("M32").Value = -307.3333
$ws.Range("N32").Value = -1052

# row 129
$ws.Range("H129").Value = 838.07935
$ws.Range("J129").Value = 963.48
$ws.Range("L129").Value = 2890.44
$ws.Range("N129").Value = -12890.44

# row 132
$ws.Range("H132").Value = 37044350
$ws.Range("I132").Value = 41673644
$ws.Range("J132").Value = 10002
$ws.Range("K132").Value = 125020932
$ws.Range("L132").Value = 30006
$ws.Range("M132").Value = -125018402
$ws.Range("N132").Value = -35066

# row 137
$ws.Range("H137").Value = 3271.9143
$ws.Range("I137").Value = 1528.0625
$ws.Range("J137").Value = 4740.421
$ws.Range("K137").Value = 4584.1875
$ws.Range("L137").Value = 14221.263
$ws.Range("M137").Value = -2034.1875
$ws.Range("N137").Value = -19321.263

# row 141
$ws.Range("H141").Value = 37639.07
$ws.Range("I141").Value = 44921.26
$ws.Range("K141").Value = 134763.78
$ws.Range("M141").Value = -129583.78

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4495.377
$ws.Range("I32").Value = 4175.08
$ws.Range("J32").Value = 5951.273
$ws.Range("K32").Value = 4175.08
$ws.Range("L32").Value = 5951.273
$ws.Range("M32").Value = -3888.08
$ws.Range("N32").Value = -6525.273

# row 132
$ws.Range("H132").Value = 2427.9778
$ws.Range("I132").Value = 1446
$ws.Range("J132").Value = 4845.154
$ws.Range("K132").Value = 4338
$ws.Range("L132").Value = 14535.462
$ws.Range("M132").Value = -1808
$ws.Range("N132").Value = -19595.462

$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 1833.5316
$ws.Range("I134").Value = 1136.8103
$ws.Range("K134").Value = 3410.4309
$ws.Range("M134").Value = -875.4309000000003

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 6581341
$ws.Range("I31").Value = 1271.0204
$ws.Range("J31").Value = 18522950
$ws.Range("K31").Value = 1271.0204
$ws.Range("L31").Value = 18522950
$ws.Range("M31").Value = -976.0204000000001
$ws.Range("N31").Value = -18523540

# row 34
$ws.Range("H34").Value = 6581341
$ws.Range("I34").Value = 1271.0204
$ws.Range("J34").Value = 18522950
$ws.Range("K34").Value = 1271.0204
$ws.Range("L34").Value = 18522950
$ws.Range("M34").Value = -1069.0204
$ws.Range("N34").Value = -18523354

# row 39
$ws.Range("H39").Value = 25255.273
$ws.Range("J39").Value = 25255.273
$ws.Range("L39").Value = 25255.273
$ws.Range("N39").Value = -26037.273

# row 41
$ws.Range("H41").Value = 28313.77
$ws.Range("I41").Value = 8500
$ws.Range("J41").Value = 29964.916
$ws.Range("K41").Value = 8500
$ws.Range("L41").Value = 29964.916
$ws.Range("M41").Value = -8072
$ws.Range("N41").Value = -30820.916

# row 42
$ws.Range("H42").Value = 33172.4
$ws.Range("J42").Value = 33172.4
$ws.Range("L42").Value = 33172.4
$ws.Range("N42").Value = -34358.4

# row 49
$ws.Range("H49").Value = 25255.273
$ws.Range("J49").Value = 25255.273
$ws.Range("L49").Value = 25255.273
$ws.Range("N49").Value = -25619.273

# row 54
$ws.Range("H54").Value = 11974
$ws.Range("J54").Value = 11974
$ws.Range("L54").Value = 11974
$ws.Range("N54").Value = -13290

# row 55
$ws.Range("H55").Value = 40063
$ws.Range("I55").Value = 40063
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 40063
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -39748
$ws.Range("N55").ClearContents()

# row 58
$ws.Range("H58").Value = 1688.2446
$ws.Range("I58").Value = 1591.8383
$ws.Range("J58").Value = 1940.3846
$ws.Range("K58").Value = 1591.8383
$ws.Range("L58").Value = 1940.3846
$ws.Range("M58").Value = -1388.8383
$ws.Range("N58").Value = -2346.3846

# row 132
$ws.Range("H132").Value = 2644.2122
$ws.Range("I132").Value = 2214.2173
$ws.Range("J132").Value = 3633.2
$ws.Range("K132").Value = 6642.651899999999
$ws.Range("L132").Value = 10899.6
$ws.Range("M132").Value = -4112.651899999999
$ws.Range("N132").Value = -15959.6

# row 134
$ws.Range("H134").Value = 4464.641
$ws.Range("I134").Value = 6545.6113
$ws.Range("J134").Value = 2680.9524
$ws.Range("K134").Value = 19636.8339
$ws.Range("L134").Value = 8042.8572
$ws.Range("M134").Value = -17101.8339
$ws.Range("N134").Value = -13112.8572

# row 136
$ws.Range("H136").Value = 1688.2446
$ws.Range("I136").Value = 1591.8383
$ws.Range("J136").Value = 1940.3846
$ws.Range("K136").Value = 4775.5149
$ws.Range("L136").Value = 5821.1538
$ws.Range("M136").Value = -2225.5149
$ws.Range("N136").Value = -10921.1538

$ws = $wb.Worksheets.Item("CUL")
# row 96
$ws.Range("H96").Value = 6967
$ws.Range("J96").Value = 6967
$ws.Range("L96").Value = 20901
$ws.Range("N96").Value = -25019

# row 131
$ws.Range("H131").Value = 824.2125
$ws.Range("J131").Value = 871.8
$ws.Range("L131").Value = 2615.4
$ws.Range("N131").Value = -12695.4

# row 137
$ws.Range("H137").Value = 2524.1904
$ws.Range("J137").Value = 3116.6667
$ws.Range("L137").Value = 9350.000100000001
$ws.Range("N137").Value = -19550.0001

$ws = $wb.Worksheets.Item("GSM")
# row 123
$ws.Range("H123").Value = 10890.833
$ws.Range("J123").Value = 10890.833
$ws.Range("L123").Value = 10890.833
$ws.Range("N123").Value = -15790.833

$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 4901.222
$ws.Range("I40").Value = 3229.8262
$ws.Range("J40").Value = 7858.3076
$ws.Range("K40").Value = 3229.8262
$ws.Range("L40").Value = 7858.3076
$ws.Range("M40").Value = -3093.8262
$ws.Range("N40").Value = -8130.3076

# row 55
$ws.Range("H55").Value = 477.72223
$ws.Range("I55").Value = 270.9
$ws.Range("J55").Value = 736.25
$ws.Range("K55").Value = 270.9
$ws.Range("L55").Value = 736.25
$ws.Range("M55").Value = -97.89999999999998
$ws.Range("N55").Value = -1082.25

# row 93
$ws.Range("H93").Value = 4832657.5
$ws.Range("I93").Value = 11112477
$ws.Range("J93").Value = 2027.4615
$ws.Range("K93").Value = 11112477
$ws.Range("L93").Value = 2027.4615
$ws.Range("M93").Value = -11111229
$ws.Range("N93").Value = -4523.461499999999

# row 132
$ws.Range("H132").Value = 4006.0667
$ws.Range("I132").Value = 1419.0286
$ws.Range("K132").Value = 4257.085800000001
$ws.Range("M132").Value = -1727.085800000001

# row 136
$ws.Range("H136").Value = 2920.75
$ws.Range("I136").Value = 1708.5172
$ws.Range("J136").Value = 7942.857
$ws.Range("K136").Value = 5125.5516
$ws.Range("L136").Value = 23828.571
$ws.Range("M136").Value = -2575.5516
$ws.Range("N136").Value = -28928.571

$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 2400.9395
$ws.Range("I122").Value = 1534.0869
$ws.Range("J122").Value = 4394.7
$ws.Range("K122").Value = 4602.2607
$ws.Range("L122").Value = 13184.1
$ws.Range("M122").Value = -2152.2607
$ws.Range("N122").Value = -18084.1

# row 132
$ws.Range("H132").Value = 6537532
$ws.Range("I132").Value = 694.34283
$ws.Range("J132").Value = 20836866
$ws.Range("K132").Value = 2083.02849
$ws.Range("L132").Value = 62510598
$ws.Range("M132").Value = 446.9715099999999
$ws.Range("N132").Value = -62515658

# row 136
$ws.Range("H136").Value = 1378.909
$ws.Range("I136").Value = 805.3714
$ws.Range("J136").Value = 2382.6
$ws.Range("K136").Value = 2416.1142
$ws.Range("L136").Value = 7147.799999999999
$ws.Range("M136").Value = 133.8858
$ws.Range("N136").Value = -12247.8
